# Insert a new data row at row 481 (pushing the existing rows 481-513 down
# to 482-514), then populate the newly inserted row with the new record.
# This matches the diff: dimension grows from A1:R513 to A1:R514 and every
# row from 481 to 513 shifts down by one, with a brand-new row of data
# appearing at row 481.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 481:513 down to 482:514, carrying formatting (incl. the date
# number format used by column D) down with them.
$ws.Rows("481:481").Insert()

# Populate the newly-inserted row 481 with the new record's data.
$ws.Cells.Item(481, 1).Value = 10
$ws.Cells.Item(481, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(481, 3).Value = "La Araucanía"
$ws.Cells.Item(481, 4).Value = 44714
$ws.Cells.Item(481, 5).Value = 9
$ws.Cells.Item(481, 6).Value = 100112043
$ws.Cells.Item(481, 7).Value = "Pepino ensalada"
$ws.Cells.Item(481, 8).Value = "Sin especificar"
$ws.Cells.Item(481, 9).Value = "Primera"
$ws.Cells.Item(481, 10).Value = 155
$ws.Cells.Item(481, 11).Value = 20000
$ws.Cells.Item(481, 12).Value = 20000
$ws.Cells.Item(481, 13).Value = 20000
$ws.Cells.Item(481, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(481, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(481, 16).Value = 333
$ws.Cells.Item(481, 17).Value = 60
$ws.Cells.Item(481, 18).Value = "Hortaliza"
